# Refresh the cryptos price/volume snapshot (GitHub Actions scheduled update).
# Price cells in column D are plain display text (e.g. "24.962.46",
# "1.706.65") that look numeric to Excel's auto-detection, so for those we
# force the cell to Text format before writing the value and then restore
# the default "Normal" style afterwards (matching the workbook's original
# un-styled price cells) to avoid Excel silently coercing them to floats
# (which would mangle values like "1.001" -> 1 or "20.60" -> 20.6).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "24.962.46"
$ws.Range("E2").Value = "  +2.21%  "
$ws.Range("D3").Value = "1.706.65"
$ws.Range("E3").Value = "  +1.36%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "315.99"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.10%  "
$ws.Range("E6").Value = "  +0.10%  "
$ws.Range("E7").Value = "  +1.77%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4027"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.57%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.483"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.08%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "52.68"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.89%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.001"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.19%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08818"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.81%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "26.16"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.04%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.471"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.01%  "
$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.997"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.18%  "
$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001356"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.16%  "
$ws.Range("D17").Value = "1.715.28"
$ws.Range("E17").Value = "  +2.67%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "96.26"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.42%  "
$ws.Range("E19").Value = "  -0.18%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "20.60"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +4.72%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.356"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.48%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9999"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.08%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "14.45"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.23%  "
$ws.Range("D24").Value = "24.954.62"
$ws.Range("E24").Value = "  +2.25%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.967"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.57%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.353"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.55%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "23.65"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +5.21%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.234"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +16.26%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "161.73"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.59%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "150.51"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +9.09%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.427"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.96%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.613"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +32.62%  "
$ws.Range("D33").Value = "1.906.73"
$ws.Range("E33").Value = "  +2.80%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08573"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.96%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.046"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.22%  "
$ws.Range("B36").Value = "VeChain"
$ws.Range("C36").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.03129"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +5.06%  "
$ws.Range("B37").Value = "InternetComputer(DFINITY)"
$ws.Range("C37").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "7.203"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.73%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2857"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.88%  "
$ws.Range("B39").Value = "Stellar"
$ws.Range("C39").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.09540"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.45%  "
$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "10.84"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.78%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8258"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.58%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "14.01"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.14%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.490"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.39%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "17.36"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.26%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.691"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.75%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.7389"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.70%  "
$ws.Range("E47").Value = "  -0.21%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.417"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.67%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.08757"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +8.83%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.001"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.15%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "139.19"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.01%  "
